{"js": "// Corresponding Author Information currently reads \"Dave Bridges PhD \".\n// Add a comma after the name: \"Dave Bridges, PhD \".\n// Find the unique run of text \"Dave Bridges PhD\" (the author-list \"Dave\n// Bridges\" mention has no trailing \"PhD\", so this search string is unique)\n// and then, scoped to that match, find the \" PhD\" sub-string and insert a\n// comma immediately before it. This only inserts a single \",\" character and\n// leaves all surrounding text/formatting untouched.\nconst body = context.document.body;\n\nconst matches = body.search(\"Dave Bridges PhD\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"Dave Bridges PhD\" in the document body.');\n}\n\nconst fullMatch = matches.items[0];\n\nconst innerMatches = fullMatch.search(\" PhD\", { matchCase: true, matchWholeWord: false });\ninnerMatches.load(\"items\");\nawait context.sync();\n\nif (innerMatches.items.length === 0) {\n  throw new Error('Could not find \" PhD\" within the \"Dave Bridges PhD\" match.');\n}\n\nconst phdRange = innerMatches.items[0];\nconst insertionPoint = phdRange.getRange(\"Start\");\ninsertionPoint.insertText(\",\", \"Start\");\nawait context.sync();\n", "ps1": "# Corresponding Author Information currently reads \"Dave Bridges PhD \".\n# Add a comma after the name: \"Dave Bridges, PhD \".\n#\n# \"Dave Bridges PhD\" (unique in the document - the author-list mention of\n# \"Dave Bridges\" has no trailing \"PhD\") is located with Find, then a second,\n# narrower Find for \" PhD\" inside that match pinpoints the exact insertion\n# point. Collapsing the inner match range to its start and calling\n# InsertBefore adds a single \",\" character without retyping/replacing any\n# surrounding text or formatting.\n\n$doc = $word.ActiveDocument\n\n$outer = $doc.Content\n$outerFind = $outer.Find\n$outerFind.Text = \"Dave Bridges PhD\"\n$outerFind.MatchCase = $true\n$outerFind.Execute()\n\nif (-not $outerFind.Found) {\n    throw 'Could not find \"Dave Bridges PhD\" in the document.'\n}\n\n$inner = $outer.Duplicate\n$innerFind = $inner.Find\n$innerFind.Text = \" PhD\"\n$innerFind.MatchCase = $true\n$innerFind.Execute()\n\nif (-not $innerFind.Found) {\n    throw 'Could not find \" PhD\" within the \"Dave Bridges PhD\" match.'\n}\n\n$inner.Collapse(1)  # wdCollapseStart\n$inner.InsertBefore(\",\")\n"}
